$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, [string]$val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
$ws.Range('D2').Value = '44.271.07'
$ws.Range('E2').Value = '  +2.09%  '

# Row 3
$ws.Range('D3').Value = '2.365.42'
$ws.Range('E3').Value = '  +0.34%  '

# Row 4
$ws.Range('E4').Value = '  -0.01%  '

# Row 5
Set-TextCell 'D5' '0.682'
$ws.Range('E5').Value = '  +4.78%  '

# Row 6
Set-TextCell 'D6' '241.53'
$ws.Range('E6').Value = '  +3.38%  '

# Row 7
Set-TextCell 'D7' '74.61'
$ws.Range('E7').Value = '  +8.28%  '

# Row 8
$ws.Range('E8').Value = '  +0.00%  '

# Row 9
Set-TextCell 'D9' '0.567'
$ws.Range('E9').Value = '  +23.41%  '

# Row 10
Set-TextCell 'D10' '0.103'
$ws.Range('E10').Value = '  +6.72%  '

# Row 11
Set-TextCell 'D11' '31.69'
$ws.Range('E11').Value = '  +19.98%  '

# Row 12
$ws.Range('E12').Value = '  +1.98%  '

# Row 13
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell 'D13' '7.20'
$ws.Range('E13').Value = '  +15.18%  '

# Row 14
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '2.718.21'
$ws.Range('E14').Value = '  +0.21%  '

# Row 15
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell 'D15' '17.13'
$ws.Range('E15').Value = '  +8.68%  '

# Row 16
Set-TextCell 'D16' '0.923'
$ws.Range('E16').Value = '  +8.63%  '

# Row 17
$ws.Range('D17').Value = '2.367.98'
$ws.Range('E17').Value = '  +0.28%  '

# Row 18
$ws.Range('D18').Value = '44.273.13'
$ws.Range('E18').Value = '  +2.05%  '

# Row 19
$ws.Range('E19').Value = '  +4.84%  '

# Row 20
Set-TextCell 'D20' '78.70'
$ws.Range('E20').Value = '  +6.47%  '

# Row 21
Set-TextCell 'D21' '6.63'
$ws.Range('E21').Value = '  +5.03%  '

# Row 22
Set-TextCell 'D22' '257.73'
$ws.Range('E22').Value = '  +2.92%  '

# Row 23
Set-TextCell 'D23' '1.00'
$ws.Range('E23').Value = '  +0.06%  '

# Row 24
Set-TextCell 'D24' '3.79'
$ws.Range('E24').Value = '  -6.79%  '

# Row 25
Set-TextCell 'D25' '2.55'
$ws.Range('E25').Value = '  +3.37%  '

# Row 26
Set-TextCell 'D26' '10.63'
$ws.Range('E26').Value = '  +6.78%  '

# Row 27
Set-TextCell 'D27' '2.31'
$ws.Range('E27').Value = '  +1.51%  '

# Row 28
Set-TextCell 'D28' '22.78'
$ws.Range('E28').Value = '  +1.21%  '

# Row 29
$ws.Range('B29').Value = 'ImmutableX'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell 'D29' '1.61'
$ws.Range('E29').Value = '  +2.61%  '

# Row 30
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 'D30' '174.49'
$ws.Range('E30').Value = '  +1.33%  '

# Row 31
Set-TextCell 'D31' '0.133'
$ws.Range('E31').Value = '  +3.55%  '

# Row 32
Set-TextCell 'D32' '0.134'
$ws.Range('E32').Value = '  +5.43%  '

# Row 33
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell 'D33' '5.29'
$ws.Range('E33').Value = '  +5.16%  '

# Row 34
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell 'D34' '0.0746'
$ws.Range('E34').Value = '  +7.98%  '

# Row 35
Set-TextCell 'D35' '5.27'
$ws.Range('E35').Value = '  +4.41%  '

# Row 36
Set-TextCell 'D36' '3.86'
$ws.Range('E36').Value = '  +5.48%  '

# Row 37
Set-TextCell 'D37' '2.47'
$ws.Range('E37').Value = '  -1.16%  '

# Row 38
Set-TextCell 'D38' '6.58'
$ws.Range('E38').Value = '  -0.11%  '

# Row 39
Set-TextCell 'D39' '0.0274'
$ws.Range('E39').Value = '  +7.58%  '

# Row 40
Set-TextCell 'D40' '19.45'
$ws.Range('E40').Value = '  +6.49%  '

# Row 41
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 'D41' '8.97'
$ws.Range('E41').Value = '  -1.21%  '

# Row 42
$ws.Range('B42').Value = 'BinanceUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextCell 'D42' '1.00'
$ws.Range('E42').Value = '  -0.11%  '

# Row 43
Set-TextCell 'D43' '1.27'
$ws.Range('E43').Value = '  +4.21%  '

# Row 44
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell 'D44' '0.190'
$ws.Range('E44').Value = '  +14.58%  '

# Row 45
$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell 'D45' '0.0999'
$ws.Range('E45').Value = '  +4.43%  '

# Row 46
$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell 'D46' '1.18'
$ws.Range('E46').Value = '  +0.77%  '

# Row 47
Set-TextCell 'D47' '101.12'
$ws.Range('E47').Value = '  +1.66%  '

# Row 48
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell 'D48' '2.44'
$ws.Range('E48').Value = '  +7.30%  '

# Row 49
$ws.Range('B49').Value = 'FTXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextCell 'D49' '4.50'
$ws.Range('E49').Value = '  +0.97%  '

# Row 50
$ws.Range('D50').Value = '1.452.20'
$ws.Range('E50').Value = '  +0.30%  '

# Row 51
$ws.Range('B51').Value = 'TerraClassic'
$ws.Range('C51').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
Set-TextCell 'D51' '0.000206'
$ws.Range('E51').Value = '  +2.69%  '
